# Update "Horarios" workbook for Linea 141 with new scrape data (02:51:07)

$wb = $excel.ActiveWorkbook

$oldTime = "02:26:12"
$newTime = "02:51:07"

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 3"

$ws1.Range("A6").Value = $newTime
$ws1.Range("D6").Value = 7

$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "03:48"
$ws1.Range("D7").Value = 57

$ws1.Range("A8").Value = $newTime
$ws1.Range("B8").Value = "04:46"
$ws1.Range("C8").Value = "215A_EL PATO"
$ws1.Range("D8").Value = 115
$ws1.Range("E8").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A6").Value = $newTime
$ws2.Range("D6").Value = 7

$ws2.Range("A7").Value = $newTime
$ws2.Range("B7").Value = "04:46"
$ws2.Range("C7").Value = "215A_EL PATO"
$ws2.Range("D7").Value = 115
$ws2.Range("E7").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
